$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit swaps the (species-observation) data between row 3 and row 4,
# while a handful of columns (C, D, K, N, P, T, U, V, W, AD, AE, AF, AG, AT,
# AW, AX, AY) stay identical between the two rows and therefore need no
# change. Columns that differ: A, B, E, F, G, H, I, J, L, Q, R, S, Y, AA.
#
# Helper: write a value that must stay a TEXT cell even though it "looks"
# like a number/date (columns I, Y, AA). Plain `Range.Value = "50"` gets
# auto-coerced to a number by Excel; going through a `="text"` formula and
# then baking it down to a literal via Copy/PasteSpecial(values) keeps the
# cell as plain text without leaving a stray NumberFormat/quote-prefix style
# behind.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($rng, [string]$text)
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# --- Row 3 gets what used to be row 4's data ---
$ws.Range("A3").Value = 111350031
$ws.Range("B3").Value = 96251
$ws.Range("E3").Value = 220093
$ws.Range("F3").Value = "Korallrot"
$ws.Range("G3").Value = "Corallorhiza trifida"
$ws.Range("H3").Value = "Châtel."
Set-TextValue $ws.Range("I3") "50"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("Q3").Value = 523281.2550886287
$ws.Range("R3").Value = 6619910.01413854
$ws.Range("S3").Value = 25
Set-TextValue $ws.Range("Y3") "2022-06-20"
Set-TextValue $ws.Range("AA3") "2022-06-20"

# --- Row 4 gets what used to be row 3's data ---
$ws.Range("A4").Value = 111350516
$ws.Range("B4").Value = 89369
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = "Vedticka"
$ws.Range("G4").Value = "Fuscoporia viticola"
$ws.Range("H4").Value = "(Schwein.) Murrill"
Set-TextValue $ws.Range("I4") "1"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("Q4").Value = 523319.7728949333
$ws.Range("R4").Value = 6619811.373445455
$ws.Range("S4").Value = 10
Set-TextValue $ws.Range("Y4") "2023-08-03"
Set-TextValue $ws.Range("AA4") "2023-08-03"

# --- Column L: the blank placeholder cell moves from row 4 to row 3 ---
$ws.Range("L4").Value = ""

Write-Host "done"
